# "tx crescimento longo prazo" -- update long-term market growth rate forecast
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width widened (21.1640625 -> 26) ---
$ws.Columns("B").ColumnWidth = 26

# --- Update forecast growth rates (2022-2055) with the new long-term figures ---
# Rows 14 & 15 keep their original values but pick up the "0.000" number format
# (style index 1), matching the rest of the historical/forecast series.
$ws.Range("B14").NumberFormat = "0.000"
$ws.Range("B14").Value = 0.03

$ws.Range("B15").NumberFormat = "0.000"
$ws.Range("B15").Value = 0.029

$ws.Range("B16").NumberFormat = "0.000"
$ws.Range("B16").Value = 0.02745165197166899

$ws.Range("B17").NumberFormat = "0.000"
$ws.Range("B17").Value = 0.029352266959988826

$ws.Range("B18").NumberFormat = "0.000"
$ws.Range("B18").Value = 0.03135258199435631

$ws.Range("B19").NumberFormat = "0.000"
$ws.Range("B19").Value = 0.0332985344113117

$ws.Range("B20").NumberFormat = "0.000"
$ws.Range("B20").Value = 0.03523025006804725

$ws.Range("B21").NumberFormat = "0.000"
$ws.Range("B21").Value = 0.03724407264518037

$ws.Range("B22").NumberFormat = "0.000"
$ws.Range("B22").Value = 0.039052382137174435

$ws.Range("B23").NumberFormat = "0.000"
$ws.Range("B23").Value = 0.040576022724664984

$ws.Range("B24").NumberFormat = "0.000"
$ws.Range("B24").Value = 0.04066469703362374

$ws.Range("B25").NumberFormat = "0.000"
$ws.Range("B25").Value = 0.04081193112810855

$ws.Range("B26").NumberFormat = "0.000"
$ws.Range("B26").Value = 0.04067613705996043

$ws.Range("B27").NumberFormat = "0.000"
$ws.Range("B27").Value = 0.04040932346010484

$ws.Range("B28").NumberFormat = "0.000"
$ws.Range("B28").Value = 0.0401511497865108

$ws.Range("B29").NumberFormat = "0.000"
$ws.Range("B29").Value = 0.03992794287336143

$ws.Range("B30").NumberFormat = "0.000"
$ws.Range("B30").Value = 0.03978586583595467

$ws.Range("B31").NumberFormat = "0.000"
$ws.Range("B31").Value = 0.03963102492861048

$ws.Range("B32").NumberFormat = "0.000"
$ws.Range("B32").Value = 0.039398816979326456

$ws.Range("B33").NumberFormat = "0.000"
$ws.Range("B33").Value = 0.03912181943017412

$ws.Range("B34").NumberFormat = "0.000"
$ws.Range("B34").Value = 0.03882042296782995

$ws.Range("B35").NumberFormat = "0.000"
$ws.Range("B35").Value = 0.03852950672201949

$ws.Range("B36").NumberFormat = "0.000"
$ws.Range("B36").Value = 0.03828933800295964

$ws.Range("B37").NumberFormat = "0.000"
$ws.Range("B37").Value = 0.038134138778175064

$ws.Range("B38").NumberFormat = "0.000"
$ws.Range("B38").Value = 0.03807365584980826

$ws.Range("B39").NumberFormat = "0.000"
$ws.Range("B39").Value = 0.038071616838173394

$ws.Range("B40").NumberFormat = "0.000"
$ws.Range("B40").Value = 0.03809108328019506

$ws.Range("B41").NumberFormat = "0.000"
$ws.Range("B41").Value = 0.03814660093649391

$ws.Range("B42").NumberFormat = "0.000"
$ws.Range("B42").Value = 0.03826524358333527

$ws.Range("B43").NumberFormat = "0.000"
$ws.Range("B43").Value = 0.03844256789548406

$ws.Range("B44").NumberFormat = "0.000"
$ws.Range("B44").Value = 0.038669283676193045

$ws.Range("B45").NumberFormat = "0.000"
$ws.Range("B45").Value = 0.0389248727613658

$ws.Range("B46").NumberFormat = "0.000"
$ws.Range("B46").Value = 0.039163729435646166

$ws.Range("B47").NumberFormat = "0.000"
$ws.Range("B47").Value = 0.03934430428858349

# --- Remove the now-unused helper column C (formerly blank spacer column); this
#     shifts the formatted-but-empty D/E helper cells left into C/D. ---
$ws.Columns("C").Delete()

# --- Refresh the view: drop the stale scroll position and move the selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select() | Out-Null
